$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the floating-point timestamp stored for the 17:00 reading (row 13)
$ws.Range("A13").Value = 45864.70865358796

# Append the new automated reading captured at 18:00:25 (row 14)
$ws.Range("A14").Value = 45864.75029034903
$ws.Range("B14").Value = 2025
$ws.Range("C14").Value = 30
$ws.Range("D14").Value = 16.94
$ws.Range("E14").Value = 78.08
$ws.Range("F14").Value = 14.06
$ws.Range("G14").Value = 6.89
$ws.Range("H14").Value = "E"
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = "18:00:25"

# Keep the new date cell using the same date/time format as the column above it
$ws.Range("A14").NumberFormat = $ws.Range("A13").NumberFormat
